$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Local" score in both the left (C) and right (F) tables.
# Downstream formulas (Num Improvement, % Improvement) recalc automatically.
$ws.Range("C3").Value = 692746622
$ws.Range("F3").Value = 692746622

# Update the active selection to match the author's final cursor position.
$ws.Range("E14").Select()
